# Natmi following Dr Hou advice
# The LR-pairs result table gained a third interacting cluster ("ECs"),
# so the Fgf2 -> Cd44 sending/target cluster grid grows from a 2x2 (4 rows)
# to a 3x3 (9 rows) cross-join, and every row's statistics are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
# E..T   : the various expression / specificity statistics (numeric)
$rows = @(
    @{ Row=2;  A="ECs";  B="Fgf2"; C="Cd44"; D="ECs";  E=2; F=0.6666666666666666; G=0.6462393333333333; H=1.938718;          I=0.03461850536298827; J=0.03461850536298827; K=3; L=1; M=261.380203;         N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=168.9141681332513;   R=1520.227513199262;   S=0.02412451884247334;  T=0.02412451884247334 },
    @{ Row=3;  A="ECs";  B="Fgf2"; C="Cd44"; D="FAPs"; E=2; F=0.6666666666666666; G=0.6462393333333333; H=1.938718;          I=0.03461850536298827; J=0.03461850536298827; K=3; L=1; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=20.679289017768;      R=186.113601159912;    S=0.002953440218019759; T=0.002953440218019759 },
    @{ Row=4;  A="ECs";  B="Fgf2"; C="Cd44"; D="sCs";  E=2; F=0.6666666666666666; G=0.6462393333333333; H=1.938718;          I=0.03461850536298827; J=0.03461850536298827; K=3; L=1; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=52.79711957254733;    R=475.174076152926;    S=0.007540546302495172; T=0.007540546302495171 },
    @{ Row=5;  A="FAPs"; B="Fgf2"; C="Cd44"; D="ECs";  E=3; F=1;                  G=15.322826;          H=45.968478;         I=0.8208310864042159;  J=0.8208310864042158;  K=3; L=1; M=261.380203;         N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=4005.083370413678;    R=36045.7503337231;    S=0.5720106862735175;   T=0.5720106862735174 },
    @{ Row=6;  A="FAPs"; B="Fgf2"; C="Cd44"; D="FAPs"; E=3; F=1;                  G=15.322826;          H=45.968478;         I=0.8208310864042159;  J=0.8208310864042158;  K=3; L=1; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=490.321667343528;     R=4412.895006091752;   S=0.07002831339388013;  T=0.07002831339388015 },
    @{ Row=7;  A="FAPs"; B="Fgf2"; C="Cd44"; D="sCs";  E=3; F=1;                  G=15.322826;          H=45.968478;         I=0.8208310864042159;  J=0.8208310864042158;  K=3; L=1; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=1251.859852507694;    R=11266.73867256925;   S=0.1787920867368182;   T=0.1787920867368182 },
    @{ Row=8;  A="sCs";  B="Fgf2"; C="Cd44"; D="ECs";  E=3; F=1;                  G=2.698388;           H=8.095164;          I=0.1445504082327959;  J=0.1445504082327959;  K=3; L=1; M=261.380203;         N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=705.305203212764;     R=6347.746828914876;   S=0.1007325131612291;   T=0.1007325131612291 },
    @{ Row=9;  A="sCs";  B="Fgf2"; C="Cd44"; D="FAPs"; E=3; F=1;                  G=2.698388;           H=8.095164;          I=0.1445504082327959;  J=0.1445504082327959;  K=3; L=1; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=86.34687252206399;    R=777.121852698576;    S=0.01233216121636344;  T=0.01233216121636345 },
    @{ Row=10; A="sCs";  B="Fgf2"; C="Cd44"; D="sCs";  E=3; F=1;                  G=2.698388;           H=8.095164;          I=0.1445504082327959;  J=0.1445504082327959;  K=3; L=1; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=220.455652481372;     R=1984.100872332348;   S=0.03148573385520329;  T=0.0314857338552033 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$($r.Row)").Value2 = $r[$col]
    }
}
